$d = $word.ActiveDocument

# This document has Track Changes turned on; turn it off so our edits land
# as plain text (matching the target revision, which carries no w:ins/w:del
# markup) instead of being recorded as tracked insertions/deletions.
$d.TrackRevisions = $false

# --- 1) "Release 2.22.0 updated the " -> "The " ------------------------------
$d.Content.Find.Execute("Release 2.22.0 updated the ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The ", 1) | Out-Null

# --- 2) "keystore" -> "keystore has been updated in this Release" -----------
#        (inserts " has been updated in this Release" right after "keystore")
$d.Content.Find.Execute("keystore", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "keystore has been updated in this Release", 1) | Out-Null

# --- 3) "(CLU), before running any commands, update" -> "(CLU) update" ------
$d.Content.Find.Execute("(CLU), before running any commands, update", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(CLU) update", 1) | Out-Null

# --- 4) "from GitHub master if you did not update it after Release 2.22.0"
#        -> "from GitHub master before running any commands" ----------------
$d.Content.Find.Execute("from GitHub master if you did not update it after Release 2.22.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "from GitHub master before running any commands", 1) | Out-Null
